# Generate Report for Handoff
# Adds a new localization-status row (for file
# "eedced51-98f9-43e5-815b-21991076e346.md") to each of the three sheets:
#   Overview (row 3), zh-cn (row 3), de-de (row 3)

$wb = $excel.ActiveWorkbook

$commitHash = "2f71f7d7c8356ede91e9ca8d46b178724adabf91"
$newFile = "eedced51-98f9-43e5-815b-21991076e346.md"
$newFileDisplayPath = "e2e\" + $newFile
$newFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitHash + "/e2e/" + $newFile

$hyperlinkColor = 15570276  # BGR-encoded 0x6495ED == OOXML color FF6495ED

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = $newFileDisplayPath
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-29 14:45:28"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, "", "", $newFileDisplayPath)
$wsOverview.Range("B3").Font.Underline = $true
$wsOverview.Range("B3").Font.Color = $hyperlinkColor

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhHandoffFile = "eedced51-98f9-43e5-815b-21991076e346.13b466b81a035275d7a66f82bc7b10eafc94deb3.zh-cn.xlf"

$wsZhCn.Range("A3").Value = $newFile
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = $zhHandoffFile
$wsZhCn.Range("H3").Value = "2016-08-29 14:45:22"
$wsZhCn.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newFileUrl, "", "", $newFile)
$wsZhCn.Range("A3").Font.Underline = $true
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deHandoffFile = "eedced51-98f9-43e5-815b-21991076e346.13b466b81a035275d7a66f82bc7b10eafc94deb3.de-de.xlf"

$wsDeDe.Range("A3").Value = $newFile
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = $deHandoffFile
$wsDeDe.Range("H3").Value = "2016-08-29 14:45:28"
$wsDeDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newFileUrl, "", "", $newFile)
$wsDeDe.Range("A3").Font.Underline = $true
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))
